$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new row for "Securities and commodity contracts and exchanges" / "5231.5232" ---
# becomes row 14, pushing old row14 ("Other financial investment activities") down to 15
$ws.Rows.Item(14).Insert()
$ws.Range("A14").Value = "Securities and commodity contracts and exchanges"
$ws.Range("B14").Value = "5231.5232"

# --- Insert new row for "Other finance and insurance" / "521.522.524" ---
# old row16 ("Funds, trusts, and other financial vehicles") stays at 16; new row becomes 17
$ws.Rows.Item(17).Insert()
$ws.Range("A17").Value = "Other finance and insurance"
$ws.Range("B17").Value = "521.522.524"

# --- Insert new row for "Real estate" / "531" ---
# old row18 ("Real estate and rental and leasing") stays at 18; new row becomes 19
$ws.Rows.Item(19).Insert()
$ws.Range("A19").Value = "Real estate"
$ws.Range("B19").Value = "531"

# --- Insert 3 new rows for the lessor breakdown, after old "Lessors of residential..." row (now row 20) ---
# Column A values were entered first for all three rows, then column B values, which is
# reflected in the shared-string insertion order.
$ws.Rows.Item(21).Insert()
$ws.Rows.Item(22).Insert()
$ws.Rows.Item(23).Insert()

$ws.Range("A21").Value = "Lessors of nonresidential buildings (except mini-warehouses)"
$ws.Range("A22").Value = "lessors of miniwarehouses and self-storage units"
$ws.Range("A23").Value = "lessors of other real estate activities"

$ws.Range("B21").Value = "53112"
$ws.Range("B22").Value = "53113"
$ws.Range("B23").Value = "53119"

# --- Update the former "Other real estate activities" row (now row 24): B changes from 53139 to "5312.5313" ---
$ws.Range("B24").Value = "5312.5313"

# --- Append a brand-new trailing row 34: "Nature of business not allocable" (no NAICS code) ---
$ws.Range("A34").Value = "Nature of business not allocable"
$ws.Range("A34").NumberFormat = "@"
